# "added lifetimes to ships"
# Update the ship lifetime figures on Sheet2 (column D, rows 4:43)
# from 700,000 to 7,000,000, and switch the active sheet/selection
# from ship_var back to Sheet2.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("D4:D43").Value = 7000000

# Make Sheet2 the active sheet (was ship_var) and set its selection,
# matching the new workbook/sheet view state.
$ws2.Activate() | Out-Null
$ws2.Range("T12").Select() | Out-Null
